$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Item description")
$ws.Columns.Item(3).Delete()
